# museos_transformacion.xlsx — raw "Responses" sheet is renamed to
# "datos_brutos" and a fresh empty "Hoja1" sheet is added for the next
# transformation step; the old auto-fit ("bestFit") column widths on the
# response-text columns (H:M) are replaced with explicit manual widths,
# and the sheet view is reset to the top-left (no scrolled topLeftCell).

$wb = $excel.ActiveWorkbook

# --- rename the raw-data sheet -------------------------------------------
$raw = $wb.Worksheets.Item(1)
$raw.Name = "datos_brutos"

# --- add the new (empty) working sheet right after it --------------------
$hoja1 = $wb.Worksheets.Add($null, $raw)
$hoja1.Name = "Hoja1"

# make sure the raw-data sheet stays the selected/active tab and the view
# is scrolled back to the top-left corner (clears topLeftCell="L1")
$raw.Activate()

# --- columns H:M no longer auto-fit; set explicit widths -----------------
$raw.Columns.Item(8).ColumnWidth = 11.666666666666666
$raw.Columns.Item(9).ColumnWidth = 10.833333333333334
$raw.Columns.Item(10).ColumnWidth = 8.833333333333334
$raw.Columns.Item(11).ColumnWidth = 13.333333333333334
$raw.Columns.Item(12).ColumnWidth = 11.5
$raw.Columns.Item(13).ColumnWidth = 13.333333333333334
